$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "adding pizza items, buffing panama and bowler"
# L12: fedora  -> bowler
# L33: sherlock -> fedora
# L41: top hat -> sherlock  (top hat item removed from the pool)
# L43: bowler  -> panama
# L53: panama  -> top       (new "top" item added)

$ws.Range("B6").Value2  = "bowler"
$ws.Range("B15").Value2 = "fedora"
$ws.Range("B17").Value2 = "sherlock"
$ws.Range("B19").Value2 = "panama"
$ws.Range("B23").Value2 = "top"

$ws.Range("B23").Select()
